# project_management.xlsx — "added project management part to final report"
#
# Adds sub-total formulas (Expected Time / Effective Time) in new columns G
# and H for each of the six task groups in the table, and updates the
# window/view state (selection + scroll position) to reflect where the user
# was working when they made the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New sub-total formulas in columns G/H, one pair per task group.
#   Group 1: rows 2-5   -> totals written on row 5
#   Group 2: rows 6-8   -> totals written on row 8
#   Group 3: rows 9-14  -> totals written on row 14
#   Group 4: rows 15-21 -> totals written on row 21
#   Group 5: rows 22-30 -> totals written on row 30
#   Group 6: rows 31-38 -> totals written on row 38
# ---------------------------------------------------------------------
$groups = @(
    @{ Start = 2;  End = 5  },
    @{ Start = 6;  End = 8  },
    @{ Start = 9;  End = 14 },
    @{ Start = 15; End = 21 },
    @{ Start = 22; End = 30 },
    @{ Start = 31; End = 38 }
)

foreach ($g in $groups) {
    $start = $g.Start
    $end = $g.End

    $ws.Range("G$end").Formula = "=SUM(D$start`:D$end)"
    $ws.Range("H$end").Formula = "=SUM(E$start`:E$end)"
}

# ---------------------------------------------------------------------
# Window / view state: scroll the sheet up so row 6 is at the top and
# leave the selection on the last new total cell (G40, the grand total
# row), matching where editing finished.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1

$ws.Range("G40").Select()
